$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: Day 47 / Level 9
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Day 47"
$ws.Range("C12").Value = 0.10347222222222223
$ws.Range("C12").NumberFormat = "h:mm"
$ws.Range("D12").Value = "70 Days"
$ws.Range("D12").NumberFormat = "h:mm"
$ws.Range("E12").Value = "Mark Sparling"
$ws.Range("F12").Value = "Level 9"
$ws.Range("G12").Value = "https://www.reddit.com/r/gamedev/comments/362c72/i_made_a_bunch_of_music_feel_free_to_use_it_in/"

# Row 13: Day 48 / Level 8
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "Day 48"
$ws.Range("C13").Value = 0.11388888888888889
$ws.Range("C13").NumberFormat = "h:mm"
$ws.Range("D13").Value = "70 Days"
$ws.Range("E13").Value = "Mark Sparling"
$ws.Range("F13").Value = "Level 8"
$ws.Range("G13").Value = "https://www.reddit.com/r/gamedev/comments/362c72/i_made_a_bunch_of_music_feel_free_to_use_it_in/"

# Update the active selection to B15
$null = $ws.Range("B15").Select()
